# Update countries & provincias Spain
# - COVID-19 daily data refresh for the "Pais" sheet.
# - Oman's case count overtakes Nigeria's, so the two rows swap places
#   (the sheet is kept sorted descending by "Casos totales").
# - Several other countries get refreshed totals.
# - The "Datos actualizados ..." timestamp banner in A1 is bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Columns: A=Pais B=Casos totales C=Nuevos casos D=Casos activos
#          E=Recuperados F=Casos criticos G=Muertes hoy H=Muertes

# --- Header banner: refresh timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 12:35"

# --- Alemania (row 11) ---
$ws.Range("B11").Value = 180808
$ws.Range("C11").Value = 19
$ws.Range("E11").Value = 10376
$ws.Range("G11").Value = 4
$ws.Range("H11").Value = 8432

# --- India (row 13) ---
$ws.Range("B13").Value = 146371
$ws.Range("C13").Value = 1421
$ws.Range("D13").Value = 61151
$ws.Range("E13").Value = 81033

# --- Suiza (row 32) ---
$ws.Range("B32").Value = 30761
$ws.Range("C32").Value = 15
$ws.Range("E32").Value = 648

# --- Oman overtakes Nigeria: row 58 becomes Oman, row 59 becomes Nigeria ---
$ws.Range("A58").Value = "Oman"
$ws.Range("B58").Value = 8118
$ws.Range("C58").Value = 348
$ws.Range("D58").Value = 2067
$ws.Range("E58").Value = 6014
$ws.Range("H58").Value = 37

$ws.Range("A59").Value = "Nigeria"
$ws.Range("B59").Value = 8068
$ws.Range("D59").Value = 2311
$ws.Range("E59").Value = 5524
$ws.Range("H59").Value = 233

# --- Marruecos (row 61) ---
$ws.Range("B61").Value = 7556
$ws.Range("C61").Value = 24
$ws.Range("D61").Value = 4841
$ws.Range("E61").Value = 2513
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = 202

# --- Senegal (row 78) ---
$ws.Range("B78").Value = 3161
$ws.Range("C78").Value = 31
$ws.Range("D78").Value = 1565
$ws.Range("E78").Value = 1560
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 36

# --- San Marino (row 129) ---
$ws.Range("D129").Value = 275
$ws.Range("E129").Value = 349

# --- Botsuana (row 187) ---
$ws.Range("D187").Value = 20
$ws.Range("E187").Value = 14
